$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.179.75"
$ws.Range("E2").Value = "  +1.17%  "

$ws.Range("D3").Value = "3.725.30"
$ws.Range("E3").Value = "  -0.38%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "612.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "189.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.37%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.637"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("E8").Value = "  -0.24%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.722"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.11%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.161"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.77%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "58.77"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +8.36%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000291"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.48%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.68"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.58%  "

$ws.Range("D14").Value = "4.308.81"
$ws.Range("E14").Value = "  -0.72%  "

$ws.Range("D15").Value = "3.718.54"
$ws.Range("E15").Value = "  -1.10%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.36"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.92%  "

$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.127"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.25%  "

$ws.Range("E18").Value = "  -0.99%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.58%  "

$ws.Range("D20").Value = "68.898.94"
$ws.Range("E20").Value = "  +0.72%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "412.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.60"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.26%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "89.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.13%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.81%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.87%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.16%  "

$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.86%  "

$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.89%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.68"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.15%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.61"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.55%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.77"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.57%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "46.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.56%  "

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.123"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.94%  "

$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "638.56"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.45%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "65.65"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.63%  "

$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").Value = "0.0₃0832"
$ws.Range("E37").Value = "  -10.38%  "

$ws.Range("B38").Value = "TheGraph"
$ws.Range("C38").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.413"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.92%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.24%  "

$ws.Range("E40").Value = "  -0.10%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.141"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.55%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.05"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.26%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0445"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.30%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.43%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.140"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.60%  "

$ws.Range("D46").Value = "2.867.75"
$ws.Range("E46").Value = "  +4.23%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.15"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.50%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.74"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.26%  "

$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "142.37"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.35%  "

$ws.Range("B50").Value = "ApeXProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.08"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.02%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.59"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -20.04%  "
